$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Incident Reports")
$ws.Range("AB1").Value = "test - -51:44"
$ws.Range("AB1").ClearContents()
Write-Host "done"
